# Purchase Test suites updated
# - Login sheet credentials changed (username/password rows)
# - Active selection on the Login sheet moved

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")
$ws.Activate() | Out-Null

# Update the stored test credentials (row 2: Username / Password columns)
$ws.Range("A2").Value = "priti"
$ws.Range("B2").Value = "jibe@123"

# Move / collapse the active selection to E6
$ws.Range("E6").Select() | Out-Null
